$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet 7: 保險 (Insurance) - add header row + 7 new columns
#   (property_category, category, date, legislator_name,
#    legislator_id, source_file, index) to rows 1-8
# ============================================================
$ws7 = $wb.Worksheets.Item(7)

# Apply the existing header/data formatting to the new columns
# (E:K) before writing values, so the added cells match the
# look of the existing B:D columns.
$ws7.Range("B1").Copy() | Out-Null
$ws7.Range("E1:K1").PasteSpecial(-4122) | Out-Null
$ws7.Range("B2").Copy() | Out-Null
$ws7.Range("E2:K8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 1
$ws7.Range("B1").Value = "company"
$ws7.Range("C1").Value = "name"
$ws7.Range("D1").Value = "owner"
$ws7.Range("E1").Value = "property_category"
$ws7.Range("F1").Value = "category"
$ws7.Range("G1").Value = "date"
$ws7.Range("H1").Value = "legislator_name"
$ws7.Range("I1").Value = "legislator_id"
$ws7.Range("J1").Value = "source_file"
$ws7.Range("K1").Value = "index"

# Row 2
$ws7.Range("A2").Value = 109
$ws7.Range("B2").Value = "全球人壽"
$ws7.Range("C2").Value = "理財型增額終身壽險"
$ws7.Range("D2").Value = "蔡瓊姿"
$ws7.Range("E2").Value = "insurance"
$ws7.Range("F2").Value = "normal"
$ws7.Range("G2").Value = "2013-12-11"
$ws7.Range("H2").Value = "吳育仁"
$ws7.Range("I2").Value = 1734
$ws7.Range("J2").Value = "tmpbcc11"
$ws7.Range("K2").Value = 109

# Row 3
$ws7.Range("A3").Value = 110
$ws7.Range("B3").Value = "全球人壽"
$ws7.Range("C3").Value = "防癌終身健康保險"
$ws7.Range("D3").Value = "蔡瓊姿"
$ws7.Range("E3").Value = "insurance"
$ws7.Range("F3").Value = "normal"
$ws7.Range("G3").Value = "2013-12-11"
$ws7.Range("H3").Value = "吳育仁"
$ws7.Range("I3").Value = 1734
$ws7.Range("J3").Value = "tmpbcc11"
$ws7.Range("K3").Value = 110

# Row 4
$ws7.Range("A4").Value = 111
$ws7.Range("B4").Value = "宏泰人壽"
$ws7.Range("C4").Value = "壽險+終身醫療險"
$ws7.Range("D4").Value = "蔡瓊姿"
$ws7.Range("E4").Value = "insurance"
$ws7.Range("F4").Value = "normal"
$ws7.Range("G4").Value = "2013-12-11"
$ws7.Range("H4").Value = "吳育仁"
$ws7.Range("I4").Value = 1734
$ws7.Range("J4").Value = "tmpbcc11"
$ws7.Range("K4").Value = 111

# Row 5
$ws7.Range("A5").Value = 112
$ws7.Range("B5").Value = "宏泰人壽"
$ws7.Range("C5").Value = "终身保險+健康保險"
$ws7.Range("D5").Value = "蔡瓊姿"
$ws7.Range("E5").Value = "insurance"
$ws7.Range("F5").Value = "normal"
$ws7.Range("G5").Value = "2013-12-11"
$ws7.Range("H5").Value = "吳育仁"
$ws7.Range("I5").Value = 1734
$ws7.Range("J5").Value = "tmpbcc11"
$ws7.Range("K5").Value = 112

# Row 6
$ws7.Range("A6").Value = 113
$ws7.Range("B6").Value = "宏泰人壽"
$ws7.Range("C6").Value = "終身壽險理財型"
$ws7.Range("D6").Value = "蔡瓊姿"
$ws7.Range("E6").Value = "insurance"
$ws7.Range("F6").Value = "normal"
$ws7.Range("G6").Value = "2013-12-11"
$ws7.Range("H6").Value = "吳育仁"
$ws7.Range("I6").Value = 1734
$ws7.Range("J6").Value = "tmpbcc11"
$ws7.Range("K6").Value = 113

# Row 7
$ws7.Range("A7").Value = 114
$ws7.Range("B7").Value = "富邦人壽"
$ws7.Range("C7").Value = "外幣終身保險"
$ws7.Range("D7").Value = "蔡瓊姿"
$ws7.Range("E7").Value = "insurance"
$ws7.Range("F7").Value = "normal"
$ws7.Range("G7").Value = "2013-12-11"
$ws7.Range("H7").Value = "吳育仁"
$ws7.Range("I7").Value = 1734
$ws7.Range("J7").Value = "tmpbcc11"
$ws7.Range("K7").Value = 114

# Row 8
$ws7.Range("A8").Value = 115
$ws7.Range("B8").Value = "南山人壽"
$ws7.Range("C8").Value = "意外骨折及特定手術傷害保險"
$ws7.Range("D8").Value = "吳育仁"
$ws7.Range("E8").Value = "insurance"
$ws7.Range("F8").Value = "normal"
$ws7.Range("G8").Value = "2013-12-11"
$ws7.Range("H8").Value = "吳育仁"
$ws7.Range("I8").Value = 1734
$ws7.Range("J8").Value = "tmpbcc11"
$ws7.Range("K8").Value = 115

# ============================================================
# Sheet 8: 債務 (Debt) - add header row + 7 new columns
#   (property_category, category, date, legislator_name,
#    legislator_id, source_file, index) to rows 1-3
# ============================================================
$ws8 = $wb.Worksheets.Item(8)

$ws8.Range("B1").Copy() | Out-Null
$ws8.Range("H1:N1").PasteSpecial(-4122) | Out-Null
$ws8.Range("B2").Copy() | Out-Null
$ws8.Range("H2:N3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 1
$ws8.Range("B1").Value = "species"
$ws8.Range("C1").Value = "debtor"
$ws8.Range("D1").Value = "owner"
$ws8.Range("E1").Value = "total"
$ws8.Range("F1").Value = "register_date"
$ws8.Range("G1").Value = "register_reason"
$ws8.Range("H1").Value = "property_category"
$ws8.Range("I1").Value = "category"
$ws8.Range("J1").Value = "date"
$ws8.Range("K1").Value = "legislator_name"
$ws8.Range("L1").Value = "legislator_id"
$ws8.Range("M1").Value = "source_file"
$ws8.Range("N1").Value = "index"

# Row 2
$ws8.Range("A2").Value = 125
$ws8.Range("B2").Value = "房屋貸款"
$ws8.Range("C2").Value = "吳育仁"
$ws8.Range("D2").Value = "玉山銀行臺北市中山區民生東路3段115號"
$ws8.Range("E2").Value = 12819300
$ws8.Range("F2").Value = "97年05月07日"
$ws8.Range("G2").Value = "購屋"
$ws8.Range("H2").Value = "debt"
$ws8.Range("I2").Value = "normal"
$ws8.Range("J2").Value = "2013-12-11"
$ws8.Range("K2").Value = "吳育仁"
$ws8.Range("L2").Value = 1734
$ws8.Range("M2").Value = "tmpbcc11"
$ws8.Range("N2").Value = 125

# Row 3
$ws8.Range("A3").Value = 126
$ws8.Range("B3").Value = "房屋貸款"
$ws8.Range("C3").Value = "蔡瓊姿"
$ws8.Range("D3").Value = "國泰人壽保險司臺北市大安區仁愛路4段296號"
$ws8.Range("E3").Value = 41033751
$ws8.Range("F3").Value = "95年11月07日"
$ws8.Range("G3").Value = "購屋"
$ws8.Range("H3").Value = "debt"
$ws8.Range("I3").Value = "normal"
$ws8.Range("J3").Value = "2013-12-11"
$ws8.Range("K3").Value = "吳育仁"
$ws8.Range("L3").Value = 1734
$ws8.Range("M3").Value = "tmpbcc11"
$ws8.Range("N3").Value = 126

Write-Output "Sheet7 and Sheet8 updated."